$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Add the new "iH2" (hydrogen) product row (row 21) of data to the four
# elasticity / production-parameter sheets that describe capital-labour(-energy)
# substitution, and refresh the view/selection state on every sheet to match
# where the workbook was left after the edits + a calculation run.
# ---------------------------------------------------------------------------

function Set-TextCellKeepStyle {
    param($ws, [string]$addr, [string]$text, [string]$styleSourceAddr)
    $ws.Range($addr).Value = $text
    $ws.Range($styleSourceAddr).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

# ---- elasKL sheet (row 21: iH2) -------------------------------------------
$wsKL = $wb.Worksheets.Item("elasKL")
Set-TextCellKeepStyle $wsKL "A21" "iH2" "A20"
$wsKL.Range("B21").Value = 0.30209999999999998

# ---- prodKL sheet (row 21: iH2) --------------------------------------------
$wsProdKL = $wb.Worksheets.Item("prodKL")
Set-TextCellKeepStyle $wsProdKL "A21" "iH2" "A20"
$wsProdKL.Range("B21").Value = 1
$wsProdKL.Range("C21").Value = 1

# ---- elasKL-E sheet (row 21: iH2) ------------------------------------------
$wsKLE = $wb.Worksheets.Item("elasKL-E")
Set-TextCellKeepStyle $wsKLE "A21" "iH2" "A20"
$wsKLE.Range("B21").Value = 0.2757
$wsKLE.Range("C21").Value = 0.30209999999999998
$wsKLE.Range("D21").Value = 0.5

# ---- prodKL-E sheet (row 21: iH2) ------------------------------------------
$wsProdKLE = $wb.Worksheets.Item("prodKL-E")
Set-TextCellKeepStyle $wsProdKLE "A21" "iH2" "A20"
$wsProdKLE.Range("B21").Value = 1
$wsProdKLE.Range("C21").Value = 1
$wsProdKLE.Range("D21").Value = 1

# ---------------------------------------------------------------------------
# Restore / update each sheet's selection (cursor position) to reflect the
# state the file was saved in, and move the active tab from elasKL-E to
# elasFU_CES.
# ---------------------------------------------------------------------------

$wsCES = $wb.Worksheets.Item("elasFU_CES")
$wsLES = $wb.Worksheets.Item("elasFU_LES")
$wsTRADE = $wb.Worksheets.Item("elasTRADE")

$wsTRADE.Activate() | Out-Null
$wsTRADE.Range("H26").Select() | Out-Null

$wsKL.Activate() | Out-Null
$wsKL.Range("B24").Select() | Out-Null

$wsProdKL.Activate() | Out-Null
$wsProdKL.Range("C25").Select() | Out-Null

$wsKLE.Activate() | Out-Null
$wsKLE.Range("B21:D21").Select() | Out-Null

$wsProdKLE.Activate() | Out-Null
$wsProdKLE.Range("I15").Select() | Out-Null

$wsLES.Activate() | Out-Null
$wsLES.Range("A2:A17").Select() | Out-Null

# elasFU_CES becomes the active / selected tab last, matching tabSelected on sheet1.
$wsCES.Activate() | Out-Null
$wsCES.Range("E27").Select() | Out-Null
